# Tracker de resultados - actualización automática
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: fill in result / profit for a match that has now resolved ---
$ws.Range("G13").Value = "Acierto"
$ws.Range("H13").Value = 0.62

# --- Row 20: fill in result / profit for a match that has now resolved ---
$ws.Range("G20").Value = "Acierto"
$ws.Range("H20").Value = 1.1

# --- Row 23: brand-new pick appended to the tracker ---
$ws.Range("A23").Value = 14655087

# Force column B to stay plain text so the date-like string isn't
# reinterpreted as an Excel date serial, then strip the formatting
# change back off so the cell keeps the sheet's default style.
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "2025-09-19"
$ws.Range("B23").ClearFormats()

$ws.Range("C23").Value = "Facundo Bagnis"
$ws.Range("D23").Value = "Federico Agustin Gomez"
$ws.Range("E23").Value = "Gana Federico Agustin Gomez"
$ws.Range("F23").Value = 1.8
# G23 / H23 stay blank (result not known yet), matching the other
# still-pending rows which carry an explicit empty text cell.
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = ""
